$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 155, pushing existing rows 155:173 down to 156:174
$ws.Rows("155:155").Insert()

# Populate the newly inserted row 155 with the new weekly price record
$ws.Range("A155").Value = 4
$ws.Range("B155").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C155").Value = "Los Lagos"
$ws.Range("D155").Value = 45223
$ws.Range("E155").Value = 10
$ws.Range("F155").Value = 100112026
$ws.Range("G155").Value = "Haba"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 120
$ws.Range("K155").Value = 16000
$ws.Range("L155").Value = 16000
$ws.Range("M155").Value = 16000
$ws.Range("N155").Value = "$/saco 25 kilos"
$ws.Range("O155").Value = "Región Metropolitana"
$ws.Range("P155").Value = 640
$ws.Range("Q155").Value = 25
$ws.Range("R155").Value = "Hortaliza"
